$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 4.9
$ws.Range("P2").Value = 2.34
$ws.Range("Q2").Value = 1.66

# Row 3
$ws.Range("K3").Value = 5.3
$ws.Range("O3").Value = 1.28
$ws.Range("T3").Value = 2.06

# Row 4
$ws.Range("K4").Value = 3.95

# Row 7
$ws.Range("G7").Value = 8.199999999999999
$ws.Range("K7").Value = 4.9

# Row 9
$ws.Range("Q9").Value = 2.5

# Row 10
$ws.Range("F10").Value = 1.18
$ws.Range("G10").Value = 1.24
$ws.Range("H10").Value = 11.5
$ws.Range("J10").Value = 7.2
$ws.Range("K10").Value = 10.5
$ws.Range("N10").Value = 1.1
$ws.Range("O10").Value = 1.09
$ws.Range("P10").Value = 3.4
$ws.Range("Q10").Value = 1.27
$ws.Range("R10").Value = 2.06
$ws.Range("S10").Value = 1.68
$ws.Range("T10").Value = 1.79
$ws.Range("U10").Value = 2
$ws.Range("W10").Value = 5
$ws.Range("X10").Value = 60
$ws.Range("AB10").Value = 18.5
$ws.Range("AI10").Value = 150
$ws.Range("AL10").Value = 32
$ws.Range("AN10").Value = 2.82
$ws.Range("AO10").Value = 190

# Row 11
$ws.Range("K11").Value = 4.6

# Row 12
$ws.Range("F12").Value = 1.49
$ws.Range("G12").Value = 1.5
$ws.Range("I12").Value = 7.4
$ws.Range("N12").Value = 6
$ws.Range("T12").Value = 1.75
$ws.Range("U12").Value = 2.26
$ws.Range("AA12").Value = 200
$ws.Range("AF12").Value = 10.5
$ws.Range("AL12").Value = 25
$ws.Range("AM12").Value = 85
$ws.Range("AN12").Value = 5.4

# Row 13
$ws.Range("F13").Value = 2.96
$ws.Range("G13").Value = 3.5
$ws.Range("H13").Value = 2.1
$ws.Range("I13").Value = 2.4
$ws.Range("J13").Value = 3.35
$ws.Range("K13").Value = 4.6
$ws.Range("L13").Value = 1.01
$ws.Range("M13").Value = 1.01
$ws.Range("N13").Value = 1.01
$ws.Range("O13").Value = 1.16
$ws.Range("P13").Value = 2.28
$ws.Range("Q13").Value = 1.45
$ws.Range("R13").Value = 1.51
$ws.Range("S13").Value = 2.04
$ws.Range("T13").Value = 1.01
$ws.Range("U13").Value = 1.01
$ws.Range("V13").Value = 1.71
$ws.Range("W13").Value = 1.4
$ws.Range("X13").Value = 38
$ws.Range("Y13").Value = 22
$ws.Range("Z13").Value = 26
$ws.Range("AA13").Value = 44
$ws.Range("AB13").Value = 28
$ws.Range("AC13").Value = 15
$ws.Range("AD13").Value = 17.5
$ws.Range("AE13").Value = 30
$ws.Range("AF13").Value = 40
$ws.Range("AG13").Value = 21
$ws.Range("AH13").Value = 21
$ws.Range("AI13").Value = 40
$ws.Range("AJ13").Value = 75
$ws.Range("AK13").Value = 46
$ws.Range("AL13").Value = 50
$ws.Range("AM13").Value = 75
$ws.Range("AN13").Value = 1000
$ws.Range("AO13").Value = 1000

# Row 14
$ws.Range("G14").Value = 2.2
$ws.Range("H14").Value = 3.85
$ws.Range("J14").Value = 3.3
$ws.Range("L14").Value = 1.45
$ws.Range("N14").Value = 3.25
$ws.Range("Q14").Value = 2.08
$ws.Range("S14").Value = 3.95
$ws.Range("W14").Value = 1.83
$ws.Range("X14").Value = 14.5
